$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - data types
$ws.Range("B1").Value = "varchar(50)"
$ws.Range("C1").Value = "varchar(50)"
$ws.Range("D1").Value = "varchar(50)"
$ws.Range("F1").Value = "int(255)"
$ws.Range("G1").Value = "int(255)"
$ws.Range("H1").Value = "varchar(8)"
$ws.Range("I1").Value = "varchar(50)"
$ws.Range("J1").Value = "varchar(50)"

# Row 2 - field labels
$ws.Range("A2").Value = "(id)입력x"
$ws.Range("B2").Value = "공장ID"
$ws.Range("C2").Value = "법인ID"
$ws.Range("D2").Value = "version"
$ws.Range("E2").Value = "사업장ID"
$ws.Range("F2").Value = "수량"
$ws.Range("G2").Value = "금액"
$ws.Range("H2").Value = "년월"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
